$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 16-57: new periodo labels (Oct 2016 .. Mar 2020, ascending) plus
# the corresponding Valor Mora / Salario Basico values for "parte 1" of
# the new estado de cuenta that this commit adds.
$rows = @(
    @{Row=16; Periodo="1610"; F=27578; G=781242},
    @{Row=17; Periodo="1611"; F=27578; G=781242},
    @{Row=18; Periodo="1612"; F=27578; G=781242},
    @{Row=19; Periodo="1701"; F=27578; G=781242},
    @{Row=20; Periodo="1702"; F=27578; G=781242},
    @{Row=21; Periodo="1703"; F=27578; G=781242},
    @{Row=22; Periodo="1704"; F=27578; G=781242},
    @{Row=23; Periodo="1705"; F=27578; G=781242},
    @{Row=24; Periodo="1706"; F=27578; G=781242},
    @{Row=25; Periodo="1707"; F=27578; G=781242},
    @{Row=26; Periodo="1708"; F=27578; G=781242},
    @{Row=27; Periodo="1709"; F=27578; G=781242},
    @{Row=28; Periodo="1710"; F=27578; G=781242},
    @{Row=29; Periodo="1711"; F=27578; G=781242},
    @{Row=30; Periodo="1712"; F=27578; G=781242},
    @{Row=31; Periodo="1801"; F=27578; G=781242},
    @{Row=32; Periodo="1802"; F=27578; G=781242},
    @{Row=33; Periodo="1803"; F=27578; G=781242},
    @{Row=34; Periodo="1804"; F=27578; G=781242},
    @{Row=35; Periodo="1805"; F=27578; G=781242},
    @{Row=36; Periodo="1806"; F=27578; G=781242},
    @{Row=37; Periodo="1807"; F=27578; G=781242},
    @{Row=38; Periodo="1808"; F=27578; G=781242},
    @{Row=39; Periodo="1809"; F=31249; G=781242},
    @{Row=40; Periodo="1810"; F=31249; G=781242},
    @{Row=41; Periodo="1811"; F=31249; G=781242},
    @{Row=42; Periodo="1812"; F=31249; G=781242},
    @{Row=43; Periodo="1901"; F=31249; G=781242},
    @{Row=44; Periodo="1902"; F=31249; G=781242},
    @{Row=45; Periodo="1903"; F=31249; G=781242},
    @{Row=46; Periodo="1904"; F=31249; G=781242},
    @{Row=47; Periodo="1905"; F=31249; G=781242},
    @{Row=48; Periodo="1906"; F=31249; G=781242},
    @{Row=49; Periodo="1907"; F=31249; G=781242},
    @{Row=50; Periodo="1908"; F=31249; G=781242},
    @{Row=51; Periodo="1909"; F=31249; G=781242},
    @{Row=52; Periodo="1910"; F=31249; G=781242},
    @{Row=53; Periodo="1911"; F=31249; G=781242},
    @{Row=54; Periodo="1912"; F=31249; G=781242},
    @{Row=55; Periodo="2001"; F=31249; G=781242},
    @{Row=56; Periodo="2002"; F=31249; G=781242},
    @{Row=57; Periodo="2003"; F=31249; G=781242}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 5).Value = $r.Periodo   # column E - Periodo Mora
    $ws.Cells.Item($r.Row, 6).Value = $r.F         # column F - Valor Mora
    $ws.Cells.Item($r.Row, 7).Value = $r.G         # column G - Salario Basico
}
